$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly record before the current row 55, shifting the
# existing rows 55-59 down to 56-60 (dimension grows from R59 to R60).
$ws.Rows.Item(55).Insert()

$ws.Range("A55").Value = 8
$ws.Range("B55").Value = "Terminal La Palmera de La Serena"
$ws.Range("C55").Value = "Coquimbo"
$ws.Range("D55").Value = 44522
$ws.Range("E55").Value = 4
$ws.Range("F55").Value = 100112052
$ws.Range("G55").Value = "Albahaca"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Primera"
$ws.Range("J55").Value = 800
$ws.Range("K55").Value = 3000
$ws.Range("L55").Value = 3500
$ws.Range("M55").Value = 3250
$ws.Range("N55").Value = "$/paquete"
$ws.Range("O55").Value = "Región de Arica y Parinacota"
$ws.Range("P55").Value = 3250
$ws.Range("Q55").Value = 1
$ws.Range("R55").Value = "Hortaliza"
